$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark. In the source document it
#    sits right after "unser Sound fuer das Ende des Timers. " - Word
#    re-creates it at the point of the most recent edit, which by the
#    end of this script is the brand-new final (empty) paragraph.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Move to the very end of the document (right after
#    "unser Sound fuer das Ende des Timers. ") and start a new
#    paragraph there.
# ------------------------------------------------------------------
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)

# --- paragraph: "Dieser wurde durch die function audio.play() ..." + hyperlink
$r.InsertAfter("Dieser wurde durch die function audio.play() mit der davor gesetzten Variable var audio gespielt. Diese haben wir von der Quelle: ")
$r.Collapse(0)
$linkStart = $r.Start
$r.InsertAfter("http://stackoverflow.com/questions/9419263/playing-audio-with-javascript")
$linkEnd = $r.End
$linkRange = $d.Range($linkStart, $linkEnd)
$d.Hyperlinks.Add($linkRange, "http://stackoverflow.com/questions/9419263/playing-audio-with-javascript")
$linkRange.Style = "Link"
$linkRange.Font.Size = 10

# --- paragraph: "Da der Sound gespielt werden soll, ..."
$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$r2.Collapse(0)
$r2.InsertAfter("Da ")
$r2.Collapse(0)
$r2.InsertAfter("der Sound gespielt werden soll, wenn der Timer auf 0 gelaufen ist, mussten wir audio noch in calculateAndShow einbinden, da sonst der Sound nach dem Schließen von alert („Fertig!“) gespielt wurde.")

# --- two empty paragraphs (no runs at all, just paragraph marks) -------
# A placeholder character is typed and then deleted again so the
# runtime does not leave a stray empty <w:r> behind in the paragraph.
$r3 = $d.Content
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$r3.Collapse(0)
$r3.InsertAfter("X")
$tmp3 = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$tmp3.Text = ""

$r4 = $d.Content
$r4.Collapse(0)
$r4.InsertParagraphAfter()
$r4.Collapse(0)
$r4.InsertAfter("X")
$tmp4 = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$tmp4.Text = ""

# --- paragraph: "Neues Rezept : lasagne"
$r5 = $d.Content
$r5.Collapse(0)
$r5.InsertParagraphAfter()
$r5.Collapse(0)
$r5.InsertAfter("Neues Rezept : lasagne")

# --- paragraph: hyperlink to eatsmarter.de
$r6 = $d.Content
$r6.Collapse(0)
$r6.InsertParagraphAfter()
$r6.Collapse(0)
$link2Start = $r6.Start
$r6.InsertAfter("http://eatsmarter.de/rezepte/lasagne")
$link2End = $r6.End
$link2Range = $d.Range($link2Start, $link2End)
$d.Hyperlinks.Add($link2Range, "http://eatsmarter.de/rezepte/lasagne")
$link2Range.Style = "Link"
$link2Range.Font.Size = 10

# ------------------------------------------------------------------
# 3. Start a final, empty paragraph and re-create the "_GoBack"
#    bookmark there (matching Word's behaviour of tracking the last
#    edit position). A temporary character is used to work around
#    the fact that a zero-length range sitting exactly on a
#    paragraph-mark cannot be used directly with Bookmarks.Add.
# ------------------------------------------------------------------
$r7 = $d.Content
$r7.Collapse(0)
$r7.InsertParagraphAfter()
$r7.Collapse(0)
$r7.InsertAfter("X")

$tmpStart = $d.Content.End - 2
$tmpEnd = $d.Content.End - 1
$tmpRange = $d.Range($tmpStart, $tmpEnd)
$d.Bookmarks.Add("_GoBack", $tmpRange)
$tmpRange.Text = ""
